$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Septiembre de 2020 a las 21:08"

# --- Row-wise data refresh ---
# The underlying source data was re-pulled; as part of that refresh two
# countries (Jordania, Burkina Faso) moved up in the case-count ranking,
# pushing the countries that used to occupy those rows down by one spot.
# We therefore rewrite each affected row in full (country name + stats).

function Set-Row($row, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Estados Unidos (row 4) - updated totals
Set-Row 4 "Estados Unidos" 6697835 21234 3960963 2538564 0 180 198308

# Italia (row 23) - only active-cases column changed
$ws.Cells.Item(23, 3).Value = 1458

# Canada (row 29) - updated totals
Set-Row 29 "Canada" 136642 501 120430 7041 0 1 9171

# Jordania now ranks ahead of Sri Lanka and Guadalupe
Set-Row 134 "Jordania" 3314 252 2206 1084 0 2 24
Set-Row 135 "Sri Lanka" 3234 39 2996 226 0 0 12
Set-Row 136 "Guadalupe" 3080 0 837 2219 0 0 24

# Yemen (row 152) - updated totals
Set-Row 152 "Yemen" 2011 2 1212 216 0 1 583

# Burkina Faso now ranks ahead of Togo and Republica de Chipre
Set-Row 156 "Burkina Faso" 1707 193 1135 516 0 0 56
Set-Row 157 "Togo" 1555 0 1189 329 0 0 37
Set-Row 158 "Republica de Chipre" 1523 0 1281 220 0 0 22

# Republica del Chad (row 165) - updated totals
Set-Row 165 "Republica del Chad" 1084 1 938 66 0 0 80
